$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Revenue (target variable) description to reflect the new
# validation/test split notebooks wording.
$ws.Range("B19").Value2 = "Zielvariable: (0)/1 bedeutet (keine)/Transaktion am Ende eines Website Besuchs"

# The shortened/rewrapped text means most description rows now fit on a
# single line, so their row heights shrink back down to the sheet's
# default wrapped-row height (15pt). Row 7 keeps its taller, multi-line
# height.
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 15
$ws.Rows.Item(6).RowHeight = 15
$ws.Rows.Item(8).RowHeight = 15
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(10).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 15
$ws.Rows.Item(12).RowHeight = 15
$ws.Rows.Item(13).RowHeight = 15
$ws.Rows.Item(14).RowHeight = 15
$ws.Rows.Item(15).RowHeight = 15
$ws.Rows.Item(16).RowHeight = 15
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(18).RowHeight = 15
$ws.Rows.Item(19).RowHeight = 15

# Reposition the view: scrolled down so row 7 is at the top, with B19
# (the cell we just edited) as the active/selected cell.
$ws.Range("B19").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1 | Out-Null
